$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Make room for the new "--spec false / No specs" row right after the
#    existing "ng g c {name_component}" options block (old row 23 -> 24).
#    Inserting a whole row pushes everything from the old row 23 down by one,
#    which reproduces every row-shift seen in the diff (25->26, 27->28, ...,
#    39->40) in a single operation.
# ---------------------------------------------------------------------------
$ws.Rows(23).Insert()

$ws.Range("C23").Value = "    --spec false"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("E23").Value = "No specs"
$ws.Range("E23").NumberFormat = "@"
# The row-insert copies the formatting of the row above into D23 (it had no
# value in the source diff); drop it so the cell doesn't linger empty.
$ws.Range("D23").Clear()

# ---------------------------------------------------------------------------
# 2) Append the brand new "PIPE" section at the bottom of the sheet
#    (new rows 41 separator, 43 header, 45 command).
# ---------------------------------------------------------------------------

# Row 41: blank separator row styled like the other section dividers
# (row 36 here is the divider that used to be row 35 before the insert above).
$ws.Range($ws.Cells(36,3), $ws.Cells(36,7)).Copy() | Out-Null
$ws.Range($ws.Cells(41,3), $ws.Cells(41,7)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 43: bold section title, same style as the other "Creamos ..." headers.
$ws.Range("B43").Value = "Creamos un nuevo PIPE"
$ws.Range("B43").Font.Bold = $true

# Row 45: the actual command text.
$ws.Range("B45").Value = "ng g p {nombre_pipe}"

# ---------------------------------------------------------------------------
# 3) Update the view so it matches the new scroll/selection position.
# ---------------------------------------------------------------------------
$ws.Range("B48").Select() | Out-Null
